$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New region row: "Eureka" added under "Region of residence" header (A2)
$ws.Range("A2").Value = "Eureka"

# Selection left where the author's cursor ended up after the edit
[void]$ws.Range("A3").Select()
